# Add a new worksheet ("Sheet2") right after the existing "Sheet1",
# populate it with the new table of data, and leave the selection/active
# sheet the way the author left it when they saved the file.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet after Sheet1 (Excel names it "Sheet2" automatically).
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# Fill in the new table: a small 2x2 block of header + data.
$ws2.Range("A1").Value = "用量"
$ws2.Range("B1").Value = "材料"
$ws2.Range("A2").Value = 2
$ws2.Range("B2").Value = "柚子酱"

# Leave the cursor on the row below the data, matching where the author's
# selection ended up, and make Sheet2 the active (visible) tab.
$ws2.Range("A3").Select() | Out-Null
